# issue #5: property land done
#
# 1. Normalises stray whitespace / stray dashes that had crept into a
#    handful of shared strings (lot/building numbers, dates, bank branch
#    names) on the 土地, 建物, 汽車 and 存款 sheets.
# 2. Re-labels the 土地 (land) sheet's header row with the common English
#    field names (name, area, share_portion, owner, register_date,
#    register_reason, acquire_value) instead of the Chinese captions.
# 3. Extends the 土地 sheet with the English metadata columns
#    (property_category, category, date, legislator_name, legislator_id,
#    source_file, index) that the other sheets (e.g. 股票) already carry.

$wb = $excel.ActiveWorkbook

$wsLand     = $wb.Worksheets.Item(1)   # 土地
$wsBuilding = $wb.Worksheets.Item(2)   # 建物
$wsCar      = $wb.Worksheets.Item(3)   # 汽車
$wsDeposit  = $wb.Worksheets.Item(5)   # 存款

# ---------------------------------------------------------------------
# 1. Re-label the 土地 header row in English
# ---------------------------------------------------------------------
$wsLand.Range("B1").Value = "name"
$wsLand.Range("C1").Value = "area"
$wsLand.Range("D1").Value = "share_portion"
$wsLand.Range("E1").Value = "owner"
$wsLand.Range("F1").Value = "register_date"
$wsLand.Range("G1").Value = "register_reason"
$wsLand.Range("H1").Value = "acquire_value"

# ---------------------------------------------------------------------
# 2. Clean up whitespace / formatting glitches in existing string values
# ---------------------------------------------------------------------

# 土地 (Land) sheet
$wsLand.Range("B2").Value = "高雄市鼓山區龍華段八小段21900000地號"
$wsLand.Range("D2").Value = "10000分之43"
$wsLand.Range("F2").Value = "85年02月05日"
$wsLand.Range("B3").Value = "高雄市鼓山區龍華段八小段21900000地號"
$wsLand.Range("D3").Value = "10000分之43"
$wsLand.Range("F3").Value = "85年02月05日"

# 建物 (Building) sheet
$wsBuilding.Range("B2").Value = "高雄市鼓山區龍華段八小段06140000建號"
$wsBuilding.Range("F2").Value = "85年02月05日"
$wsBuilding.Range("B3").Value = "高雄市鼓山區龍華段八小段06140000建號"
$wsBuilding.Range("F3").Value = "85年02月05日"

# 汽車 (Car) sheet
$wsCar.Range("E2").Value = "98年10月12日"
$wsCar.Range("B3").Value = "納智捷L91ML(客車）"
$wsCar.Range("E3").Value = "99年05月03日"

# 存款 (Deposit) sheet
$wsDeposit.Range("B2").Value = "國泰世華商業銀行南高雄分行"
$wsDeposit.Range("B3").Value = "台北富邦商業銀行基隆路分行"
$wsDeposit.Range("B4").Value = "合作金庫商業銀行港都分行"
$wsDeposit.Range("B5").Value = "台北公館郵局（第13支局）"
$wsDeposit.Range("B7").Value = "高雄金福路郵局(第44支局）"

# ---------------------------------------------------------------------
# 3. Extend the 土地 (Land) sheet with the shared metadata columns
#    (I:O), matching the layout already used on the 股票 sheet.
# ---------------------------------------------------------------------

# Force the "date" column to plain text first so "2012-04-06" isn't
# auto-converted into a date serial number by the COM layer.
$wsLand.Range("K1:K3").NumberFormat = "@"

$wsLand.Range("I1").Value = "property_category"
$wsLand.Range("J1").Value = "category"
$wsLand.Range("K1").Value = "date"
$wsLand.Range("L1").Value = "legislator_name"
$wsLand.Range("M1").Value = "legislator_id"
$wsLand.Range("N1").Value = "source_file"
$wsLand.Range("O1").Value = "index"

$wsLand.Range("I2").Value = "land"
$wsLand.Range("J2").Value = "normal"
$wsLand.Range("K2").Value = "2012-04-06"
$wsLand.Range("L2").Value = "林國正"
$wsLand.Range("M2").Value = 1742
$wsLand.Range("N2").Value = "tmpd6491"
$wsLand.Range("O2").Value = 13

$wsLand.Range("I3").Value = "land"
$wsLand.Range("J3").Value = "normal"
$wsLand.Range("K3").Value = "2012-04-06"
$wsLand.Range("L3").Value = "林國正"
$wsLand.Range("M3").Value = 1742
$wsLand.Range("N3").Value = "tmpd6491"
$wsLand.Range("O3").Value = 14

# Match the header / data cell formatting already used for columns B:H
# (bold, centered, bordered header row; plain data rows) by copying the
# formats across instead of re-deriving a brand-new style.
$wsLand.Range("H1").Copy()
$wsLand.Range("I1:O1").PasteSpecial(-4122)

$wsLand.Range("H2").Copy()
$wsLand.Range("I2:O2").PasteSpecial(-4122)

$wsLand.Range("H3").Copy()
$wsLand.Range("I3:O3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
